$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compute Engine")

$text = "Реализовано в computeEngine.spec.ts"

# Rows 12, 13, 14 get a new comment in column C (same formatting as the
# existing C column cells above them), and row heights grow to fit the
# wrapped text.
$ws.Range("C12").Value = $text
$ws.Range("C13").Value = $text
$ws.Range("C14").Value = $text

$ws.Range("C11").Copy()
$ws.Range("C12:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(12).RowHeight = 46
$ws.Rows.Item(13).RowHeight = 46
$ws.Rows.Item(14).RowHeight = 46

# Update the selection to reflect the last edited cell
$ws.Range("C14").Select()
